# SEMINAR PROJECT MANAGENMENT.pptx — swap the "Qualitative" / "Quantitative"
# Risk Analysis title slides and merge the "Assessment "/"Matrix" runs on the
# Risk Assessment Matrix slide.

$p = $ppt.ActivePresentation

# Slide 25: "Quantitative Risk Analysis"  ->  "Qualitative Risk Analysis"
$s25 = $p.Slides.Item(25)
$tr25 = $s25.Shapes.Item(1).TextFrame.TextRange
$tr25.Characters(1, 13).Text = "Qualitative "

# Slide 26: "The Risk Assessment Matrix" -> merge "Assessment "+"Matrix" runs
$s26 = $p.Slides.Item(26)
$tr26 = $s26.Shapes.Item(1).TextFrame.TextRange
$tr26.Characters(10, 18).Text = "Assessment Matrix"

# Slide 31: "Qualitative Risk Analysis"  ->  "Quantitative Risk Analysis"
$s31 = $p.Slides.Item(31)
$tr31 = $s31.Shapes.Item(1).TextFrame.TextRange
$tr31.Characters(1, 12).Text = "Quantitative "
